$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 20000200
$ws.Range("I6").Value = 25000036
$ws.Range("J6").Value = 850
$ws.Range("K6").Value = 75000108
$ws.Range("L6").Value = 2550
$ws.Range("M6").Value = -74999996
$ws.Range("N6").Value = -2774
$ws.Range("H53").Value = 1483.2858
$ws.Range("I53").Value = 2029.2222
$ws.Range("J53").Value = 500.6
$ws.Range("K53").Value = 2029.2222
$ws.Range("L53").Value = 500.6
$ws.Range("M53").Value = -1392.2222
$ws.Range("N53").Value = -1774.6
$ws.Range("H64").Value = 3866
$ws.Range("I64").Value = 3750
$ws.Range("J64").Value = 3957.1428
$ws.Range("K64").Value = 3750
$ws.Range("L64").Value = 3957.1428
$ws.Range("M64").Value = -3502
$ws.Range("N64").Value = -4453.1428
$ws.Range("H67").Value = 3866
$ws.Range("I67").Value = 3750
$ws.Range("J67").Value = 3957.1428
$ws.Range("K67").Value = 3750
$ws.Range("L67").Value = 3957.1428
$ws.Range("M67").Value = -2892
$ws.Range("N67").Value = -5673.1428
$ws.Range("H76").Value = 3035.1667
$ws.Range("I76").Value = 3013.0435
$ws.Range("J76").Value = 3162.375
$ws.Range("K76").Value = 3013.0435
$ws.Range("L76").Value = 3162.375
$ws.Range("M76").Value = -2698.0435
$ws.Range("N76").Value = -3792.375
$ws.Range("H79").Value = 3035.1667
$ws.Range("I79").Value = 3013.0435
$ws.Range("J79").Value = 3162.375
$ws.Range("K79").Value = 3013.0435
$ws.Range("L79").Value = 3162.375
$ws.Range("M79").Value = -1921.0435
$ws.Range("N79").Value = -5346.375
$ws.Range("H86").Value = 6600.8335
$ws.Range("J86").Value = 8431.200000000001
$ws.Range("L86").Value = 8431.200000000001
$ws.Range("N86").Value = -10677.2
$ws.Range("H89").Value = 6600.8335
$ws.Range("J89").Value = 8431.200000000001
$ws.Range("L89").Value = 42156
$ws.Range("N89").Value = -53388
$ws.Range("H92").Value = 1075.65
$ws.Range("I92").Value = 913.3125
$ws.Range("J92").Value = 1725
$ws.Range("K92").Value = 913.3125
$ws.Range("L92").Value = 1725
$ws.Range("M92").Value = 334.6875
$ws.Range("N92").Value = -4221
$ws.Range("H113").Value = 4200.087
$ws.Range("I113").Value = 3529.8
$ws.Range("J113").Value = 4715.6924
$ws.Range("K113").Value = 3529.8
$ws.Range("L113").Value = 4715.6924
$ws.Range("M113").Value = -275.8000000000002
$ws.Range("N113").Value = -11223.6924
$ws.Range("H137").Value = 2806.5715
$ws.Range("I137").Value = 2761.0588
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 8283.1764
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -5733.1764
$ws.Range("N137").Value = -14100
$ws.Range("H138").Value = 5885112
$ws.Range("I138").Value = 1820.7894
$ws.Range("J138").Value = 13337281
$ws.Range("K138").Value = 5462.3682
$ws.Range("L138").Value = 40011843
$ws.Range("M138").Value = -322.3681999999999
$ws.Range("N138").Value = -40022123
$ws.Range("H140").Value = 41745
$ws.Range("J140").Value = 41745
$ws.Range("L140").Value = 41745
$ws.Range("N140").Value = -52105

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3199.2856
$ws.Range("I63").Value = 2698.75
$ws.Range("J63").Value = 3866.6667
$ws.Range("K63").Value = 2698.75
$ws.Range("L63").Value = 3866.6667
$ws.Range("M63").Value = -2012.75
$ws.Range("N63").Value = -5238.6667
$ws.Range("H66").Value = 3199.2856
$ws.Range("I66").Value = 2698.75
$ws.Range("J66").Value = 3866.6667
$ws.Range("K66").Value = 13493.75
$ws.Range("L66").Value = 19333.3335
$ws.Range("M66").Value = -10061.75
$ws.Range("N66").Value = -26197.3335
$ws.Range("H102").Value = 2154.4443
$ws.Range("I102").Value = 1815
$ws.Range("J102").Value = 2833.3333
$ws.Range("K102").Value = 1815
$ws.Range("L102").Value = 2833.3333
$ws.Range("M102").Value = -193
$ws.Range("N102").Value = -6077.3333
$ws.Range("H103").Value = 27681
$ws.Range("J103").Value = 27681
$ws.Range("L103").Value = 27681
$ws.Range("N103").Value = -30025
$ws.Range("H137").Value = 41149.832
$ws.Range("J137").Value = 41579.8
$ws.Range("L137").Value = 41579.8
$ws.Range("N137").Value = -51779.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 28091.666
$ws.Range("J51").Value = 28091.666
$ws.Range("L51").Value = 28091.666
$ws.Range("N51").Value = -29073.666
$ws.Range("H55").Value = 29670
$ws.Range("J55").Value = 29670
$ws.Range("L55").Value = 29670
$ws.Range("N55").Value = -30216

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1383.5135
$ws.Range("I58").Value = 1406.3928
$ws.Range("J58").Value = 1312.3334
$ws.Range("K58").Value = 1406.3928
$ws.Range("L58").Value = 1312.3334
$ws.Range("M58").Value = -1203.3928
$ws.Range("N58").Value = -1718.3334
$ws.Range("H62").Value = 2950
$ws.Range("I62").Value = 2330
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 2330
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -1706
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 2950
$ws.Range("I65").Value = 2330
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 11650
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -8530
$ws.Range("N65").Value = -28740
$ws.Range("H132").Value = 2255.4614
$ws.Range("I132").Value = 1838.421
$ws.Range("J132").Value = 3387.4285
$ws.Range("K132").Value = 5515.263
$ws.Range("L132").Value = 10162.2855
$ws.Range("M132").Value = -2985.263
$ws.Range("N132").Value = -15222.2855
$ws.Range("H136").Value = 1383.5135
$ws.Range("I136").Value = 1406.3928
$ws.Range("J136").Value = 1312.3334
$ws.Range("K136").Value = 4219.178400000001
$ws.Range("L136").Value = 3937.0002
$ws.Range("M136").Value = -1669.178400000001
$ws.Range("N136").Value = -9037.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3728.5715
$ws.Range("J100").Value = 3728.5715
$ws.Range("L100").Value = 11185.7145
$ws.Range("N100").Value = -12807.7145
$ws.Range("H109").Value = 2141.0454
$ws.Range("I109").Value = 859.5714
$ws.Range("J109").Value = 2739.0667
$ws.Range("K109").Value = 2578.7142
$ws.Range("L109").Value = 8217.2001
$ws.Range("M109").Value = -1538.7142
$ws.Range("N109").Value = -10297.2001
$ws.Range("H113").Value = 675.5238000000001
$ws.Range("I113").Value = 676
$ws.Range("J113").Value = 674.8889
$ws.Range("K113").Value = 2028
$ws.Range("L113").Value = 2024.6667
$ws.Range("M113").Value = 142
$ws.Range("N113").Value = -6364.6667
$ws.Range("H115").Value = 3048
$ws.Range("I115").Value = 1014
$ws.Range("J115").Value = 3500
$ws.Range("K115").Value = 3042
$ws.Range("L115").Value = 10500
$ws.Range("M115").Value = -1867
$ws.Range("N115").Value = -12850
$ws.Range("H122").Value = 1349.5
$ws.Range("J122").Value = 1681.7273
$ws.Range("L122").Value = 15135.5457
$ws.Range("N122").Value = -20035.5457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 13024.667
$ws.Range("J57").Value = 13709.6
$ws.Range("L57").Value = 13709.6
$ws.Range("N57").Value = -15349.6
$ws.Range("H110").Value = 30542.857
$ws.Range("J110").Value = 30542.857
$ws.Range("L110").Value = 30542.857
$ws.Range("N110").Value = -38722.857
$ws.Range("H132").Value = 3365.1667
$ws.Range("I132").Value = 3341.6667
$ws.Range("J132").Value = 3388.6667
$ws.Range("K132").Value = 10025.0001
$ws.Range("L132").Value = 10166.0001
$ws.Range("M132").Value = -7495.000100000001
$ws.Range("N132").Value = -15226.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 33258.65
$ws.Range("J127").Value = 33258.65
$ws.Range("L127").Value = 33258.65
$ws.Range("N127").Value = -43178.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 992.8333
$ws.Range("I100").Value = 995.5
$ws.Range("J100").Value = 987.5
$ws.Range("K100").Value = 1991
$ws.Range("L100").Value = 1975
$ws.Range("M100").Value = -1450
$ws.Range("N100").Value = -3057
$ws.Range("H107").Value = 8032.615
$ws.Range("I107").Value = 8669.5
$ws.Range("J107").Value = 390
$ws.Range("K107").Value = 26008.5
$ws.Range("L107").Value = 1170
$ws.Range("M107").Value = -24088.5
$ws.Range("N107").Value = -5010
$ws.Range("H113").Value = 616.6667
$ws.Range("J113").Value = 866.6667
$ws.Range("L113").Value = 2600.0001
$ws.Range("N113").Value = -6940.0001
$ws.Range("H132").Value = 2816.3953
$ws.Range("I132").Value = 2890.074
$ws.Range("J132").Value = 2692.0625
$ws.Range("K132").Value = 8670.222
$ws.Range("L132").Value = 8076.1875
$ws.Range("M132").Value = -6140.222
$ws.Range("N132").Value = -13136.1875
